$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1 : Title slide (AI deck)
# ---------------------------------------------------------------
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "人工智能：从图灵测试到通用人工智能 (AGI)"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Topic: AI的发展历史与未来趋势"

# ---------------------------------------------------------------
# Slide 2 : Section header
# ---------------------------------------------------------------
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "计算智能的起源与逻辑奠基 (1950-1980)"

# ---------------------------------------------------------------
# Slide 3 : Content slide + speaker notes
# ---------------------------------------------------------------
$s = $p.Slides.Item(3)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "图灵测试与符号 AI 的诞生"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "1950年：阿兰·图灵发表《计算机器与智能》，提出著名的“图灵测试” (Turing Test)。"
$tr.Paragraphs(3).Runs(1).Text = "1956年：达特茅斯会议 (Dartmouth Workshop) 正式确立“人工智能”学科，麦卡锡、明斯基等人为学科领袖。"
$tr.Paragraphs(4).Runs(1).Text = "逻辑主义时代：基于规则的专家系统（如 MYCIN）在特定医疗诊断领域取得初步成功。"
$tr.Paragraphs(5).Runs(1).Text = "瓶颈出现：早期 AI 难以处理模糊信息，导致70年代中期进入第一个“AI 冬天”。"
$s.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = "本节介绍 AI 的前身及其早期的哲学与数学基础，重点强调对人类智能的逻辑模拟尝试。"

# ---------------------------------------------------------------
# Slide 4 : Content slide + speaker notes
# ---------------------------------------------------------------
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "联结主义与神经网络的复兴"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "1986年：Rumelhart 提出反向传播算法 (Backpropagation)，解决了多层感知器的训练难题。"
$tr.Paragraphs(3).Runs(1).Text = "统计学习方法崛起：SVM 与随机森林在90年代成为机器学习的主流工具。"
$tr.Paragraphs(4).Runs(1).Text = "GPU 计算能力的增强：为复杂的矩阵运算提供了硬件基础，神经网络的研究重心逐渐转向深度化。"
$tr.Paragraphs(5).Runs(1).Text = "循环神经网络 (RNN) 与 LSTM：在高盛等金融机构及自然语言处理中开始显露头角。"
$s.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = "介绍神经网络从被冷落到重新获得学术界关注的过程，为后来的深度学习爆发做铺垫。"

# ---------------------------------------------------------------
# Slide 5 : Section header
# ---------------------------------------------------------------
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "深度学习革命与大模型时代 (2012-Present)"

# ---------------------------------------------------------------
# Slide 6 : Content slide + speaker notes
# ---------------------------------------------------------------
$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "从 ImageNet 到 Transformer 架构"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "2012年：AlexNet 以领先第二名10.8%的优势夺得 ImageNet 冠军，开启深度卷积神经网络时代。"
$tr.Paragraphs(3).Runs(1).Text = "2017年：Google 发表《Attention is All You Need》，提出 Transformer 架构，颠覆序列建模模式。"
$tr.Paragraphs(4).Runs(1).Text = "预训练大模型 (LLMs)：GPT-3 的 1750亿参数规模展示了模型容量与涌现能力 (Emergent Abilities) 的正相关性。"
$tr.Paragraphs(5).Runs(1).Text = "推理与对齐：利用 RLHF (基于人类反馈的强化学习) 解决了模型在道德与逻辑层面的幻觉问题。"
$s.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = "分析深度学习在视觉和 NLP 领域的双重突破，揭示 Transformer 架构为何成为当今万物互联的技术底座。"

# ---------------------------------------------------------------
# Slide 7 : References
# ---------------------------------------------------------------
$s = $p.Slides.Item(7)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "https://en.wikipedia.org/wiki/Main_Page"
$tr.Paragraphs(3).Runs(1).Text = "https://www.nih.gov/"
$tr.Paragraphs(4).Runs(1).Text = "https://scholar.google.com/"
$null = $tr.InsertAfter("`rhttps://www.jstor.org/")

# ---------------------------------------------------------------
# Slide 8 : Title slide (health deck)
# ---------------------------------------------------------------
$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "代谢革命：基于现代营养学的长寿科学"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Topic: 现代营养学与代谢健康科学"

# ---------------------------------------------------------------
# Slide 9 : Section header
# ---------------------------------------------------------------
$s = $p.Slides.Item(9)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "细胞能量代谢与胰岛素平衡"

# ---------------------------------------------------------------
# Slide 10 : Content slide + speaker notes
# ---------------------------------------------------------------
$s = $p.Slides.Item(10)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "胰岛素敏感性：健康的万能钥匙"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "胰岛素抵抗 (Insulin Resistance)：不仅是糖尿病前兆，更是 2 型糖尿病、多囊卵巢综合征 (PCOS) 的核心驱动因素。"
$tr.Paragraphs(3).Runs(1).Text = "血糖波动的负面影响：餐后高血糖导致的糖基化终产物 (AGEs) 会加速血管内膜老化。"
$tr.Paragraphs(4).Runs(1).Text = "低 GI 饮食策略：通过全谷物和高纤维摄入，维持血清能量供应的平滑曲线。"
$tr.Paragraphs(5).Runs(1).Text = "动态血糖监测 (CGM)：现代医疗技术从盲目补给向实时精准控糖的转变。"
$s.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = "旨在揭示胰岛素在人体内调取和储存能量的底层逻辑，以及不当饮食对代谢系统的长期损伤。"

# ---------------------------------------------------------------
# Slide 11 : Section header
# ---------------------------------------------------------------
$s = $p.Slides.Item(11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "线粒体功能与抗炎生活方式"

# ---------------------------------------------------------------
# Slide 12 : Content slide + speaker notes
# ---------------------------------------------------------------
$s = $p.Slides.Item(12)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "线粒体：细胞的能量工厂"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "线粒体自噬 (Mitophagy)：通过断食或高强度间歇训练 (HIIT) 触发坏死线粒体的自我清理。"
$tr.Paragraphs(3).Runs(1).Text = "氧化压力与抗氧化平衡：SOD 等内源性酶在抵御超氧阴离子自由基中的核心作用。"
$tr.Paragraphs(4).Runs(1).Text = "睡眠与线粒体修复：深度睡眠期是大脑清除代谢废物（β-淀粉样蛋白）的唯一窗口期。"
$tr.Paragraphs(5).Runs(1).Text = "Omega-3s 的抗炎机制：通过调节细胞膜流动性来降低慢性系统性炎症水平。"
$s.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = "探讨如何从微观层面优化身体引擎，强调生活细节对生物学年龄的逆转作用。"

# ---------------------------------------------------------------
# Slide 13 : References
# ---------------------------------------------------------------
$s = $p.Slides.Item(13)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "https://en.wikipedia.org/wiki/Main_Page"
$tr.Paragraphs(3).Runs(1).Text = "https://www.nih.gov/"
$tr.Paragraphs(4).Runs(1).Text = "https://scholar.google.com/"
$null = $tr.InsertAfter("`rhttps://www.jstor.org/")
